$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the "PickupID/POD No" column (C) on rows that were
# refreshed with the latest NAM Prod run results.
# Several of the new values look like plain numbers ("15282511", ...).
# Typing / assigning them directly would make Excel coerce them into
# numeric cells, but the source workbook stores them as text (shared
# strings) with no special number-format / quote-prefix styling.
# To reproduce that faithfully we stage each value as a text formula
# result in an unused scratch cell, copy it, and paste-special just the
# value into the destination - this preserves the text type without
# adding any quote-prefix style to the destination cell.

function Set-TextValue {
    param($range, [string]$value)
    $scratch = $ws.Range("Z100")
    $escaped = $value.Replace('"', '""')
    $scratch.Formula = "=""$escaped"""
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

Set-TextValue $ws.Range("C2")  "15282511"
Set-TextValue $ws.Range("C3")  "15282520"
Set-TextValue $ws.Range("C11") "15282548"
Set-TextValue $ws.Range("C13") "15282569"
Set-TextValue $ws.Range("C14") "15282561"
Set-TextValue $ws.Range("C24") "15282715"
Set-TextValue $ws.Range("C28") "158941792"

# "RT00006568" is not numeric-looking, so a plain assignment already
# keeps it as text without Excel adding any quote-prefix style.
$ws.Range("C27").Value = "RT00006568"

$excel.CutCopyMode = $false
